$wb = $excel.ActiveWorkbook

# Data for 2023-09-24 update: (SheetName, CellRef, NewValue) triples
# derived from the canonical OOXML diff (column J = year 2023 running totals).
$updates = @(
    @{Sheet='Citywide Totals'; Cell='J2'; Value=5590},
    @{Sheet='Citywide Totals'; Cell='J3'; Value=5969},
    @{Sheet='Citywide Totals'; Cell='J4'; Value=1295},
    @{Sheet='Citywide Totals'; Cell='J5'; Value=458},
    @{Sheet='Citywide Totals'; Cell='J6'; Value=7561},
    @{Sheet='Citywide Totals'; Cell='J7'; Value=20873},
    @{Sheet='Logan Square'; Cell='J3'; Value=39},
    @{Sheet='Logan Square'; Cell='J6'; Value=187},
    @{Sheet='Logan Square'; Cell='J7'; Value=291},
    @{Sheet='Austin'; Cell='J3'; Value=397},
    @{Sheet='Austin'; Cell='J6'; Value=436},
    @{Sheet='South Chicago'; Cell='J3'; Value=158},
    @{Sheet='South Chicago'; Cell='J6'; Value=116},
    @{Sheet='South Chicago'; Cell='J7'; Value=425},
    @{Sheet='Garfield Park'; Cell='J2'; Value=235},
    @{Sheet='Garfield Park'; Cell='J3'; Value=315},
    @{Sheet='Garfield Park'; Cell='J6'; Value=329},
    @{Sheet='Garfield Park'; Cell='J7'; Value=959},
    @{Sheet='West Pullman'; Cell='J3'; Value=109},
    @{Sheet='West Pullman'; Cell='J7'; Value=311},
    @{Sheet='Grand Crossing'; Cell='J2'; Value=191},
    @{Sheet='Grand Crossing'; Cell='J7'; Value=647},
    @{Sheet='New City'; Cell='J4'; Value=24},
    @{Sheet='New City'; Cell='J6'; Value=186},
    @{Sheet='New City'; Cell='J7'; Value=530},
    @{Sheet='Woodlawn'; Cell='J2'; Value=91},
    @{Sheet='Woodlawn'; Cell='J3'; Value=129},
    @{Sheet='Woodlawn'; Cell='J7'; Value=325},
    @{Sheet='By Neighborhood'; Cell='J4'; Value=90},
    @{Sheet='By Neighborhood'; Cell='J5'; Value=65},
    @{Sheet='By Neighborhood'; Cell='J7'; Value=612},
    @{Sheet='By Neighborhood'; Cell='J10'; Value=142},
    @{Sheet='By Neighborhood'; Cell='J15'; Value=228},
    @{Sheet='By Neighborhood'; Cell='J19'; Value=609},
    @{Sheet='By Neighborhood'; Cell='J20'; Value=432},
    @{Sheet='By Neighborhood'; Cell='J21'; Value=59},
    @{Sheet='By Neighborhood'; Cell='J27'; Value=126},
    @{Sheet='By Neighborhood'; Cell='J29'; Value=1173},
    @{Sheet='By Neighborhood'; Cell='J31'; Value=190},
    @{Sheet='By Neighborhood'; Cell='J33'; Value=959},
    @{Sheet='By Neighborhood'; Cell='J37'; Value=647},
    @{Sheet='By Neighborhood'; Cell='J42'; Value=867},
    @{Sheet='By Neighborhood'; Cell='J45'; Value=30},
    @{Sheet='By Neighborhood'; Cell='J46'; Value=70},
    @{Sheet='By Neighborhood'; Cell='J47'; Value=159},
    @{Sheet='By Neighborhood'; Cell='J48'; Value=243},
    @{Sheet='By Neighborhood'; Cell='J52'; Value=526},
    @{Sheet='By Neighborhood'; Cell='J53'; Value=291},
    @{Sheet='By Neighborhood'; Cell='J54'; Value=406},
    @{Sheet='By Neighborhood'; Cell='J57'; Value=87},
    @{Sheet='By Neighborhood'; Cell='J60'; Value=128},
    @{Sheet='By Neighborhood'; Cell='J63'; Value=79},
    @{Sheet='By Neighborhood'; Cell='J65'; Value=530},
    @{Sheet='By Neighborhood'; Cell='J67'; Value=794},
    @{Sheet='By Neighborhood'; Cell='J72'; Value=86},
    @{Sheet='By Neighborhood'; Cell='J76'; Value=304},
    @{Sheet='By Neighborhood'; Cell='J78'; Value=257},
    @{Sheet='By Neighborhood'; Cell='J79'; Value=596},
    @{Sheet='By Neighborhood'; Cell='J80'; Value=32},
    @{Sheet='By Neighborhood'; Cell='J82'; Value=27},
    @{Sheet='By Neighborhood'; Cell='J83'; Value=425},
    @{Sheet='By Neighborhood'; Cell='J84'; Value=179},
    @{Sheet='By Neighborhood'; Cell='J85'; Value=876},
    @{Sheet='By Neighborhood'; Cell='J86'; Value=127},
    @{Sheet='By Neighborhood'; Cell='J87'; Value=72},
    @{Sheet='By Neighborhood'; Cell='J88'; Value=224},
    @{Sheet='By Neighborhood'; Cell='J89'; Value=274},
    @{Sheet='By Neighborhood'; Cell='J91'; Value=232},
    @{Sheet='By Neighborhood'; Cell='J94'; Value=209},
    @{Sheet='By Neighborhood'; Cell='J95'; Value=311},
    @{Sheet='By Neighborhood'; Cell='J97'; Value=170},
    @{Sheet='By Neighborhood'; Cell='J98'; Value=151},
    @{Sheet='By Neighborhood'; Cell='J99'; Value=325},
    @{Sheet='By Neighborhood'; Cell='J101'; Value=20873},
    @{Sheet='Gage Park'; Cell='J2'; Value=73},
    @{Sheet='Gage Park'; Cell='J6'; Value=52},
    @{Sheet='Gage Park'; Cell='J7'; Value=190},
    @{Sheet='North Lawndale'; Cell='J2'; Value=197},
    @{Sheet='North Lawndale'; Cell='J3'; Value=303},
    @{Sheet='North Lawndale'; Cell='J4'; Value=61},
    @{Sheet='North Lawndale'; Cell='J6'; Value=211},
    @{Sheet='North Lawndale'; Cell='J7'; Value=794},
    @{Sheet='South Deering'; Cell='J6'; Value=54},
    @{Sheet='South Deering'; Cell='J7'; Value=179},
    @{Sheet='Loop'; Cell='J3'; Value=80},
    @{Sheet='Loop'; Cell='J7'; Value=406},
    @{Sheet='Englewood'; Cell='J2'; Value=352},
    @{Sheet='Englewood'; Cell='J3'; Value=406},
    @{Sheet='Englewood'; Cell='J6'; Value=306},
    @{Sheet='Englewood'; Cell='J7'; Value=1173},
    @{Sheet='Lake View'; Cell='J2'; Value=38},
    @{Sheet='Lake View'; Cell='J6'; Value=123},
    @{Sheet='Lake View'; Cell='J7'; Value=243},
    @{Sheet='Chatham'; Cell='J2'; Value=152},
    @{Sheet='Chatham'; Cell='J3'; Value=178},
    @{Sheet='Chatham'; Cell='J6'; Value=226},
    @{Sheet='Chatham'; Cell='J7'; Value=609},
    @{Sheet='River North'; Cell='J3'; Value=63},
    @{Sheet='River North'; Cell='J6'; Value=169},
    @{Sheet='River North'; Cell='J7'; Value=304},
    @{Sheet='Humboldt Park'; Cell='J2'; Value=190},
    @{Sheet='Humboldt Park'; Cell='J3'; Value=174},
    @{Sheet='Humboldt Park'; Cell='J6'; Value=446},
    @{Sheet='Humboldt Park'; Cell='J7'; Value=867},
    @{Sheet='Avondale'; Cell='J2'; Value=31},
    @{Sheet='Avondale'; Cell='J6'; Value=78},
    @{Sheet='Avondale'; Cell='J7'; Value=142},
    @{Sheet='Rogers Park'; Cell='J6'; Value=72},
    @{Sheet='Rogers Park'; Cell='J7'; Value=257},
    @{Sheet='Jefferson Park'; Cell='J2'; Value=21},
    @{Sheet='Jefferson Park'; Cell='J6'; Value=28},
    @{Sheet='Jefferson Park'; Cell='J7'; Value=70},
    @{Sheet='Washington Park'; Cell='J6'; Value=54},
    @{Sheet='Washington Park'; Cell='J7'; Value=232},
    @{Sheet='Chinatown'; Cell='J6'; Value=40},
    @{Sheet='Chinatown'; Cell='J7'; Value=59},
    @{Sheet='Roseland'; Cell='J3'; Value=209},
    @{Sheet='Roseland'; Cell='J6'; Value=169},
    @{Sheet='Roseland'; Cell='J7'; Value=596},
    @{Sheet='Chicago Lawn'; Cell='J3'; Value=151},
    @{Sheet='Chicago Lawn'; Cell='J6'; Value=113},
    @{Sheet='Chicago Lawn'; Cell='J7'; Value=432},
    @{Sheet='Auburn Gresham'; Cell='J2'; Value=188},
    @{Sheet='Auburn Gresham'; Cell='J3'; Value=184},
    @{Sheet='Auburn Gresham'; Cell='J5'; Value=17},
    @{Sheet='Auburn Gresham'; Cell='J7'; Value=612},
    @{Sheet='West Loop'; Cell='J2'; Value=37},
    @{Sheet='West Loop'; Cell='J6'; Value=114},
    @{Sheet='West Loop'; Cell='J7'; Value=209},
    @{Sheet='Kenwood'; Cell='J6'; Value=75},
    @{Sheet='Kenwood'; Cell='J7'; Value=159},
    @{Sheet='Brighton Park'; Cell='J2'; Value=66},
    @{Sheet='Brighton Park'; Cell='J6'; Value=95},
    @{Sheet='Brighton Park'; Cell='J7'; Value=228},
    @{Sheet='Wicker Park'; Cell='J2'; Value=26},
    @{Sheet='Wicker Park'; Cell='J6'; Value=94},
    @{Sheet='Wicker Park'; Cell='J7'; Value=151},
    @{Sheet='West Town'; Cell='J6'; Value=117},
    @{Sheet='West Town'; Cell='J7'; Value=170},
    @{Sheet='United Center'; Cell='J6'; Value=103},
    @{Sheet='United Center'; Cell='J7'; Value=224},
    @{Sheet='Uptown'; Cell='J2'; Value=86},
    @{Sheet='Uptown'; Cell='J3'; Value=75},
    @{Sheet='Uptown'; Cell='J6'; Value=82},
    @{Sheet='Uptown'; Cell='J7'; Value=274},
    @{Sheet='Armour Square'; Cell='J6'; Value=30},
    @{Sheet='Armour Square'; Cell='J7'; Value=65},
    @{Sheet='Edgewater'; Cell='J6'; Value=44},
    @{Sheet='Edgewater'; Cell='J7'; Value=126},
    @{Sheet='Streeterville'; Cell='J6'; Value=23},
    @{Sheet='Streeterville'; Cell='J7'; Value=127},
    @{Sheet='Mckinley Park'; Cell='J6'; Value=33},
    @{Sheet='Mckinley Park'; Cell='J7'; Value=87},
    @{Sheet='Morgan Park'; Cell='J2'; Value=46},
    @{Sheet='Morgan Park'; Cell='J7'; Value=128},
    @{Sheet='South Shore'; Cell='J2'; Value=229},
    @{Sheet='South Shore'; Cell='J6'; Value=255},
    @{Sheet='South Shore'; Cell='J7'; Value=876},
    @{Sheet='Old Town'; Cell='J3'; Value=26},
    @{Sheet='Old Town'; Cell='J4'; Value=8},
    @{Sheet='Old Town'; Cell='J7'; Value=86},
    @{Sheet='Sheffield & DePaul'; Cell='J5'; Value=18},
    @{Sheet='Sheffield & DePaul'; Cell='J6'; Value=27},
    @{Sheet='Jackson Park'; Cell='J6'; Value=10},
    @{Sheet='Jackson Park'; Cell='J7'; Value=30},
    @{Sheet='Rush & Division'; Cell='J3'; Value=8},
    @{Sheet='Rush & Division'; Cell='J7'; Value=32},
    @{Sheet='Little Village'; Cell='J3'; Value=162},
    @{Sheet='Little Village'; Cell='J7'; Value=526},
    @{Sheet='Archer Heights'; Cell='J6'; Value=34},
    @{Sheet='Archer Heights'; Cell='J7'; Value=90},
    @{Sheet='Ukrainian Village'; Cell='J6'; Value=47},
    @{Sheet='Ukrainian Village'; Cell='J7'; Value=72}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Host "Applied $($updates.Count) cell updates for 2023-09-24."
